$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Date: update timestamp
$ws.Range("B8").Value = "2023-09-01T17:43:23+00:00"

# Content: complete -> supplement
$ws.Range("B19").Value = "supplement"

# Supplements: add the referenced CodeSystem URL
$ws.Range("B20").Value = "http://terminology.hl7.org/CodeSystem/research-study-phase"

# Count: clear the value (was "1")
$ws.Range("B21").Value = ""
